# Update "想去人数" (want-to-go count) figures for a handful of events
# on the "展览" (exhibitions) sheet and the "全部类型" (all types) sheet,
# reflecting refreshed scrape counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 371
$ws1.Range("F4").Value  = 5225
$ws1.Range("F5").Value  = 561
$ws1.Range("F6").Value  = 10414
$ws1.Range("F9").Value  = 121
$ws1.Range("F10").Value = 133
$ws1.Range("F11").Value = 829

# --- Sheet "全部类型" ---
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F3").Value  = 371
$ws2.Range("F6").Value  = 5225
$ws2.Range("F7").Value  = 561
$ws2.Range("F9").Value  = 10414
$ws2.Range("F12").Value = 121
$ws2.Range("F15").Value = 133
$ws2.Range("F16").Value = 829
